# Case 4.77 re-run with 380 kV slack/reference voltage (vm_pu.xlsx)
# Updates the per-bus voltage magnitude results for rows 2-25 (buses 0-23):
#  - Column B (slack bus, e.g. ext_grid) drops from 1.05 p.u. to 1.02 p.u.
#  - All other bus voltage magnitudes are recomputed for the new setpoint.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ B=1.02; C=1.036050048989564; D=1.045031388801901; E=1.039697817890114; F=1.051946139633491; I=1.04126698348536; J=1.041160482676163; K=1.047800771080341; L=1.042482267786427; M=1.054696233967583; N=1.01765761782517 }
    3 = @{ B=1.02; C=1.036928480810799; D=1.045751409860684; E=1.040524407164275; F=1.05284153239199; I=1.041514128942042; J=1.041683075675493; K=1.048332388372463; L=1.043119085859744; M=1.055404161991951; N=1.017832621105605 }
    4 = @{ B=1.02; C=1.037497364333266; D=1.046217683251864; E=1.041060079568865; F=1.053421754523846; I=1.041672932125922; J=1.042021060430047; K=1.048676080161196; L=1.043531314836513; M=1.055862445333643; N=1.017945757324121 }
    5 = @{ B=1.02; C=1.037736636426826; D=1.04641379199718; E=1.041285469510949; F=1.053665880130409; I=1.041739424995225; J=1.042163107948193; K=1.04882049515871; L=1.04370465414234; M=1.056055155708937; N=1.017993294810306 }
    6 = @{ B=1.02; C=1.037776817879641; D=1.046446724599084; E=1.041323324761127; F=1.053706881587782; I=1.041750573704994; J=1.042186955899711; K=1.048844738755573; L=1.043733760779925; M=1.05608751540654; N=1.018001275086473 }
    7 = @{ B=1.02; C=1.037500561054496; D=1.046220303325; E=1.04106309048289; F=1.053425015757451; I=1.041673821660792; J=1.04202295863947; K=1.048678010129892; L=1.043533630856133; M=1.05586502015202; N=1.017946392621014 }
    8 = @{ B=1.02; C=1.036346819087422; D=1.045274645408968; E=1.039976998138279; F=1.052248566250863; I=1.041350738192311; J=1.041337129530386; K=1.047980495074492; L=1.042697448480196; M=1.054935437946691; N=1.017716781998364 }
    9 = @{ B=1.02; C=1.03431750773971; D=1.043611190165577; E=1.038069467729616; F=1.050182041332447; I=1.040772904460526; J=1.04012737339788; K=1.046749137541049; L=1.041225309383301; M=1.053299039081076; N=1.017311411273664 }
    10 = @{ B=1.02; C=1.032967218146936; D=1.042504272008851; E=1.036802112319273; F=1.048808844599126; I=1.040381995963745; J=1.039320101540812; K=1.04592679142092; L=1.040244847486702; M=1.052209294361135; N=1.017040672774942 }
    11 = @{ B=1.02; C=1.032383156920079; D=1.042025470283662; E=1.036254380460052; F=1.048215319612596; I=1.040211388362116; J=1.038970375251206; K=1.045570379211849; L=1.03982054124079; M=1.051737721909388; N=1.016923328637232 }
    12 = @{ B=1.02; C=1.032166305445523; D=1.041847698843049; E=1.036051086491558; F=1.047995021617359; I=1.040147816211972; J=1.0388404465202; K=1.045437943283174; L=1.039662972363856; M=1.051562604479924; N=1.016879725295351 }
    13 = @{ B=1.02; C=1.032212816509577; D=1.041885827918559; E=1.036094686561359; F=1.048042268855455; I=1.040161461733826; J=1.038868317779934; K=1.045466353418622; L=1.039696769712147; M=1.051600165678502; N=1.016889079100481 }
    14 = @{ B=1.02; C=1.032365229966914; D=1.042010774064488; E=1.036237572884296; F=1.04819710636462; I=1.040206137562543; J=1.038959635801847; K=1.045559433000723; L=1.039807515778925; M=1.05172324571035; N=1.016919724705706 }
    15 = @{ B=1.02; C=1.032459149537274; D=1.042087767736909; E=1.036325630851522; F=1.048292528592731; I=1.04023363721541; J=1.039015896592541; K=1.045616776002378; L=1.039875755073471; M=1.051799085460524; N=1.016938604299777 }
    16 = @{ B=1.02; C=1.033005993601996; D=1.042536059193077; E=1.036838485554216; F=1.048848257753298; I=1.040393290419172; J=1.039343308184489; K=1.045950438451697; L=1.040273012463966; M=1.052240597387555; N=1.017048458181913 }
    17 = @{ B=1.02; C=1.033349182140217; D=1.042817395842608; E=1.037160465682147; F=1.049197141794346; I=1.040493077936662; J=1.03954863946032; K=1.046159648308509; L=1.040522266736827; M=1.052517625891848; N=1.01711733676793 }
    18 = @{ B=1.02; C=1.033549418241496; D=1.042981542999544; E=1.037348371617444; F=1.04940074388101; I=1.040551152812808; J=1.039668389028794; K=1.046281644872975; L=1.040667675630617; M=1.052679240225524; N=1.017157501591905 }
    19 = @{ B=1.02; C=1.033617703721862; D=1.043037521074893; E=1.037412459662776; F=1.049470184509915; I=1.040570932859727; J=1.039709217673667; K=1.046323237085098; L=1.040717260184307; M=1.052734351294161; N=1.017171194899344 }
    20 = @{ B=1.02; C=1.033312355039336; D=1.042787206080841; E=1.03712590986807; F=1.049159699077097; I=1.040482385069755; J=1.039526611076382; K=1.046137205372972; L=1.040495521714467; M=1.052487900418581; N=1.017109947870945 }
    21 = @{ B=1.02; C=1.032320345373507; D=1.041973978405159; E=1.03619549204845; F=1.048151506047943; I=1.040192987184626; J=1.038932745596924; K=1.045532024715127; L=1.039774902781198; M=1.051687000439494; N=1.016910700791101 }
    22 = @{ B=1.02; C=1.031697179041595; D=1.041463115219655; E=1.035611416539391; F=1.047518562481113; I=1.040009869293209; J=1.038559215751639; K=1.045151243204285; L=1.039322037807414; M=1.051183707204466; N=1.016785330957339 }
    23 = @{ B=1.02; C=1.032027478795088; D=1.041733890731223; E=1.03592095875116; F=1.04785400749893; I=1.040107053465763; J=1.038757244194178; K=1.045353128883112; L=1.039562089130828; M=1.051450486952356; N=1.01685180077359 }
    24 = @{ B=1.02; C=1.033328995435651; D=1.042800847385779; E=1.037141523843056; F=1.04917661750838; I=1.040487217116225; J=1.039536564806568; K=1.046147346468512; L=1.04050760656763; M=1.052501331993159; N=1.017113286629135 }
    25 = @{ B=1.02; C=1.034841683470183; D=1.044040878288854; E=1.038561853462704; F=1.050715503844579; I=1.04092329360615; J=1.040440264393239; K=1.047067732078673; L=1.041605728247728; M=1.053721884101476; N=1.017416297730998 }
}

foreach ($rowKey in $newValues.Keys) {
    $rowVals = $newValues[$rowKey]
    foreach ($colKey in $rowVals.Keys) {
        $addr = "{0}{1}" -f $colKey, $rowKey
        $ws.Range($addr).Value = $rowVals[$colKey]
    }
}

Write-Output "case with 380 kV done"
